$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 2) ---
$ws.Range("A2").Value = "MCH163-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

# Give row 2 its own (new) font: Calibri 10pt, theme text color - matches
# the style used by the author for the newly-entered metadata row.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1

# Reuse that same formatting for the rest of the row's cells (including the
# still-empty C2/D2/H2, which the source row also carries formatting for).
$ws.Range("A2").Copy()
$ws.Range("C2:H2").PasteSpecial(-4122)

# F2 ("extentAndMedium") carries its own distinct alignment in the source.
$ws.Range("F2").HorizontalAlignment = -4108

# --- Update the active selection to match the saved workbook state ---
$ws.Range("G6").Select()
